$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 3099.503889238888
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 3747.096267775823
